# "1) Api Test examples"
#  - Fix the "TC#1_AllCounvers" typo on the AllCountriesTestData sheet
#  - Add a new "NewCountryTestData" worksheet (POST /countries example) after CountryIdTestData
#  - Leave the selections the way the saved workbook has them

$wb = $excel.ActiveWorkbook

$wsAll       = $wb.Worksheets.Item("AllCountriesTestData")
$wsCountryId = $wb.Worksheets.Item("CountryIdTestData")

# --- 1. Tweak the selection on "CountryIdTestData" ---
$wsCountryId.Activate()
$wsCountryId.Range("A2").Select()

# --- 2. Add the new "NewCountryTestData" worksheet right after "CountryIdTestData" ---
$wsNew = $wb.Worksheets.Add($null, $wsCountryId)
$wsNew.Name = "NewCountryTestData"

# Reuse the header/body formatting (fills, borders, wrap, quote-prefix, …) from the
# sibling "CountryIdTestData" sheet so the new sheet looks the same as its neighbours.
$wsCountryId.Range("A1:D2").Copy()
$wsNew.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$wsNew.Application.CutCopyMode = $false

$wsNew.Rows.Item(1).RowHeight = 22
$wsNew.Rows.Item(2).RowHeight = 128
# (values are the Excel "characters" width equivalent of the 18.6640625 / 39.33203125 /
# 31.1640625 / 27.6640625 raw column widths used by the sibling test-data sheets)
$wsNew.Columns.Item(1).ColumnWidth = 17.83
$wsNew.Columns.Item(2).ColumnWidth = 38.5
$wsNew.Columns.Item(3).ColumnWidth = 30.33
$wsNew.Columns.Item(4).ColumnWidth = 26.83

# Data row - the "create new country" POST example (entered first, ahead of the
# header, matching the order new strings were appended to the workbook). A2/C2/D2 sit
# on the quote-prefixed style (like the equivalent columns on CountryIdTestData), so a
# leading apostrophe is used to keep that format instead of falling back to plain text.
$wsNew.Range("A2").Value = "'TC#1_CreateNewCountry"
$wsNew.Range("B2").Value = "Verify creation new country using Post call, "
$wsNew.Range("C2").Value = "'name: Test Country;`n alpha2_code: TC;`nalpha3_code: TCY"
$wsNew.Range("D2").Value = "'201"

# Header row
$wsNew.Range("A1").Value = "TestCaseId"
$wsNew.Range("B1").Value = "TestCaseDescription"
$wsNew.Range("C1").Value = "inputData"
$wsNew.Range("D1").Value = "expectedStatusCode"

$wsNew.Range("C1").Select()

# --- 3. Fix the "AllCounvers" typo and re-activate "AllCountriesTestData" so it stays
#        the selected tab, matching the saved workbook state. The leading apostrophe
#        keeps the quote-prefixed cell style (A2 already uses a quotePrefix format). ---
$wsAll.Activate()
$wsAll.Range("A2").Value = "'TC#1_AllCountries"
$wsAll.Range("B9").Select()
